# Cotações atualizadas - 2025-09-18
# Adds a new row (row 14) with the latest fund quotations, matching the
# formatting/style of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date (serial 45918 = 2025-09-18), styled like the other date cells.
$ws.Range("A14").Value = 45918
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat

# New quotations for the four funds (kept as text with comma decimals,
# consistent with the rest of the sheet).
$ws.Range("B14").Value = "20,9409"
$ws.Range("C14").Value = "14,8955"
$ws.Range("D14").Value = "14,8001"
$ws.Range("E14").Value = "14,8001"
